# "removed false start data"
# The first two recorded rows (r2:r3) were a false start and are removed;
# the rows that follow (old r4:r5) shift up to become the new r2:r3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the stray rows, then delete them (shifting everything below up by 2).
$ws.Range("A2:XFD3").Select()
$ws.Rows("2:3").Delete()
